$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing formatted style from column A (e.g. A2) onto the new rows
# in column A so the new index cells match the workbook's existing style (s="1").
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A2515:A2562").PasteSpecial(-4122) | Out-Null

$newRows = @(
    @{Row=2515; A=2513; B="[0, -7, 5]"; C=15.28572233692817},
    @{Row=2516; A=2514; B="[4, -1, 7]"; C=14.23819195963161},
    @{Row=2517; A=2515; B="[3, -2, 6]"; C=14.49239652783319},
    @{Row=2518; A=2516; B="[5, 0, 7]"; C=14.2916524203619},
    @{Row=2519; A=2517; B="[2, -2, 6]"; C=15.46022536371772},
    @{Row=2520; A=2518; B="[5, 0, 6]"; C=14.38537152011593},
    @{Row=2521; A=2519; B="[3, -3, 6]"; C=14.23435105656958},
    @{Row=2522; A=2520; B="[3, -3, 5]"; C=14.33962475017112},
    @{Row=2523; A=2521; B="[3, -2, 7]"; C=14.43950852768009},
    @{Row=2524; A=2522; B="[2, -3, 6]"; C=14.68058402331847},
    @{Row=2525; A=2523; B="[2, -4, 6]"; C=14.67877272011704},
    @{Row=2526; A=2524; B="[5, -7, -2]"; C=13.98585790195392},
    @{Row=2527; A=2525; B="[4, -2, 7]"; C=14.0374308942397},
    @{Row=2528; A=2526; B="[5, -1, 7]"; C=14.20636515333507},
    @{Row=2529; A=2527; B="[2, -2, 7]"; C=15.47987207469113},
    @{Row=2530; A=2528; B="[3, -4, 5]"; C=14.34796612566152},
    @{Row=2531; A=2529; B="[3, -3, 7]"; C=14.21189433694705},
    @{Row=2532; A=2530; B="[4, -2, 6]"; C=14.09966332257641},
    @{Row=2533; A=2531; B="[4, -2, 5]"; C=14.15908284861707},
    @{Row=2534; A=2532; B="[4, -3, 5]"; C=14.06777833042299},
    @{Row=2535; A=2533; B="[3, -4, 6]"; C=14.24142405702893},
    @{Row=2536; A=2534; B="[6, 0, 7]"; C=14.68923399284561},
    @{Row=2537; A=2535; B="[6, -1, 6]"; C=14.3511995564321},
    @{Row=2538; A=2536; B="[5, -2, 5]"; C=13.98439257477286},
    @{Row=2539; A=2537; B="[6, -1, 7]"; C=14.37562758966072},
    @{Row=2540; A=2538; B="[6, -6, -1]"; C=14.35937706818516},
    @{Row=2541; A=2539; B="[6, -7, -2]"; C=14.0678058912073},
    @{Row=2542; A=2540; B="[5, -2, 4]"; C=14.15424887767322},
    @{Row=2543; A=2541; B="[7, -5, -2]"; C=14.08510830348672},
    @{Row=2544; A=2542; B="[3, 7, -2]"; C=14.43950852800636},
    @{Row=2545; A=2543; B="[4, -7, -3]"; C=13.98245047374972},
    @{Row=2546; A=2544; B="[6, -6, -3]"; C=14.01078430776325},
    @{Row=2547; A=2545; B="[7, -6, -3]"; C=14.00527186733722},
    @{Row=2548; A=2546; B="[7, -1, 6]"; C=14.36654541156535},
    @{Row=2549; A=2547; B="[7, 1, 3]"; C=14.97396957082438},
    @{Row=2550; A=2548; B="[6, -2, 3]"; C=14.15615621406356},
    @{Row=2551; A=2549; B="[7, 1, 2]"; C=14.93359506072333},
    @{Row=2552; A=2550; B="[1, 7, -2]"; C=16.20788697882952},
    @{Row=2553; A=2551; B="[1, -3, 6]"; C=15.19884109298317},
    @{Row=2554; A=2552; B="[1, -3, 7]"; C=15.21355883566591},
    @{Row=2555; A=2553; B="[2, -3, 7]"; C=14.69790549308267},
    @{Row=2556; A=2554; B="[0, -3, 7]"; C=15.33528636979073},
    @{Row=2557; A=2555; B="[2, -4, 7]"; C=14.69712175326478},
    @{Row=2558; A=2556; B="[1, -4, 7]"; C=15.1770709824043},
    @{Row=2559; A=2557; B="[1, -5, 7]"; C=15.17628761134422},
    @{Row=2560; A=2558; B="[1, -6, 7]"; C=15.17868816231102},
    @{Row=2561; A=2559; B="[1, -6, 6]"; C=15.16379453357068},
    @{Row=2562; A=2560; B="[2, -5, 7]"; C=14.69425127571442}
)

foreach ($item in $newRows) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}
